# Update countries & provincias Spain
# Applies the 26-Abril-2020 04:52 data refresh to the "Pais" sheet:
#   - bumps the "Datos actualizados" timestamp note
#   - refreshes totals for Estados Unidos (row 4) and Australia (row 46)
#   - inserts a fresh Honduras row and a fresh Guatemala row into the
#     ranking, which pushes Senegal/Uruguay and Sri Lanka/Georgia/Malta/
#     Jordania down one ranking slot each (their totals are carried over
#     unchanged to the row below)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." footer note -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 04:52"

# --- Estados Unidos (row 4) ------------------------------------------------
$ws.Range("B4").Value = 960896
$ws.Range("C4").Value = 245
$ws.Range("E4").Value = 788469
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 54265

# --- Australia (row 46) ----------------------------------------------------
$ws.Range("D46").Value = 5523
$ws.Range("E46").Value = 1104
$ws.Range("F46").Value = 42
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 83

# --- Ranking reshuffle around rows 102-110 ---------------------------------
# New Honduras entry (row 102) with fresh numbers
$ws.Range("A102").Value = "Honduras"
$ws.Range("B102").Value = 627
$ws.Range("C102").Value = 36
$ws.Range("D102").Value = 65
$ws.Range("E102").Value = 503
$ws.Range("F102").Value = 10
$ws.Range("G102").Value = 4
$ws.Range("H102").Value = 59

# Senegal shifts down from row 102 to row 103 (totals unchanged)
$ws.Range("A103").Value = "Senegal"
$ws.Range("B103").Value = 614
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 276
$ws.Range("E103").Value = 331
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 7

# Uruguay shifts down from row 103 to row 104 (totals unchanged)
$ws.Range("A104").Value = "Uruguay"
$ws.Range("B104").Value = 596
$ws.Range("C104").Value = 33
$ws.Range("D104").Value = 370
$ws.Range("E104").Value = 212
$ws.Range("F104").Value = 9
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 14

# Row 105 (San Marino) is untouched.

# New Guatemala entry (row 106) with fresh numbers
$ws.Range("A106").Value = "Guatemala"
$ws.Range("B106").Value = 473
$ws.Range("C106").Value = 43
$ws.Range("D106").Value = 45
$ws.Range("E106").Value = 415
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 13

# Sri Lanka shifts down from row 106 to row 107 (totals unchanged)
$ws.Range("A107").Value = "Sri Lanka"
$ws.Range("B107").Value = 460
$ws.Range("C107").Value = 8
$ws.Range("D107").Value = 118
$ws.Range("E107").Value = 335
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 7

# Georgia shifts down from row 107 to row 108 (totals unchanged)
$ws.Range("A108").Value = "Georgia"
$ws.Range("B108").Value = 456
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 139
$ws.Range("E108").Value = 312
$ws.Range("F108").Value = 6
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 5

# Malta shifts down from row 108 to row 109 (totals unchanged)
$ws.Range("A109").Value = "Malta"
$ws.Range("B109").Value = 448
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 249
$ws.Range("E109").Value = 195
$ws.Range("F109").Value = 2
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 4

# Jordania shifts down from row 109 to row 110 (totals unchanged)
$ws.Range("A110").Value = "Jordania"
$ws.Range("B110").Value = 444
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 332
$ws.Range("E110").Value = 105
$ws.Range("F110").Value = 5
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 7

# Row 111 (Taiwan) is untouched - the old Guatemala row is absorbed by the
# new row 106 entry above, keeping the overall row count the same.
